$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
$ws.Activate()

# --- Split "alarms(<bool>)" / "enables/disables alarms" row into two rows:
#     row 82: "Artisan Command" | (blank) | "enables/disables alarms"
#     row 83 (new): "alarm(n,<bool>)" | "enables/disables alarm number " + italic "n"

# Insert a new row before the current row 83, pushing everything below down by one.
$ws.Rows.Item(83).Insert()

# Row 82: clear the old command text in column B, keep the description in column C as is.
$ws.Range("B82").Clear()
$ws.Range("B82").NumberFormat = "General"

# Row 83 (newly inserted): set the new command name.
$ws.Range("A83").Clear()
$ws.Range("B83").Value2 = "alarm(n,<bool>)"

# Row 83, column C: rich text description with italic "n".
$prefix = "enables/disables alarm number "
$ws.Range("C83").Value2 = $prefix + "n"
$ws.Range("C83").Characters(1, $prefix.Length).Font.Italic = $false
$ws.Range("C83").Characters($prefix.Length + 1, 1).Font.Italic = $true

# --- Merge the runs of the "ramp(n,<bool>)" description (now shifted to row 109)
#     so that "toggles playback ramping per " and "event type " become a single run.
$full = "toggles playback ramping per event type n from {1,2,3,4}"
$ws.Range("C109").Value2 = $full
$italicStart = ("toggles playback ramping per event type ").Length + 1
$ws.Range("C109").Characters(1, $italicStart - 1).Font.Italic = $false
$ws.Range("C109").Characters($italicStart, 2).Font.Italic = $true
$afterItalic = $italicStart + 2
$ws.Range("C109").Characters($afterItalic, $full.Length - $afterItalic + 1).Font.Italic = $false

# --- Update the view/selection state to reflect the new row position.
$ws.Range("B83:C83").Select()
